$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articels")

# Row 4 (A4 = "fds") is emptied out
$ws.Range("A4").ClearContents()

# New row 6: A6 = "/add_vendor"
$ws.Range("A6").Value = "/add_vendor"

# New row 8: A8 = "/add_vendor"
$ws.Range("A8").Value = "/add_vendor"
